$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7 through 18 (only rows 1-6 should remain)
$ws.Range("A7:A18").EntireRow.Delete() | Out-Null

# Update the remaining values in A1:A6
$values = @(1111, 2222, 3333, 4444, 5555, 6666)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Move the selection to C12 as in the final state
$ws.Range("C12").Select() | Out-Null
